$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1540.2
$ws.Range("I100").Value = 1084.25
$ws.Range("K100").Value = 1084.25
$ws.Range("M100").Value = -543.25

$ws.Range("H107").Value = 755.6818
$ws.Range("I107").Value = 847.2105
$ws.Range("J107").Value = 176
$ws.Range("K107").Value = 847.2105
$ws.Range("L107").Value = 176
$ws.Range("M107").Value = 1072.7895
$ws.Range("N107").Value = -4016

$ws.Range("H132").Value = 944.3043
$ws.Range("I132").Value = 842.525
$ws.Range("K132").Value = 2527.575
$ws.Range("M132").Value = 2.425000000000182

$ws.Range("H137").Value = 1511.9231
$ws.Range("I137").Value = 1416.05
$ws.Range("J137").Value = 1831.5
$ws.Range("K137").Value = 4248.15
$ws.Range("L137").Value = 5494.5
$ws.Range("M137").Value = -1698.15
$ws.Range("N137").Value = -10594.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 5250
$ws.Range("J3").Value = 5500
$ws.Range("L3").Value = 5500
$ws.Range("N3").Value = -5730

$ws.Range("H22").Value = 7248.4
$ws.Range("I22").Value = 4737.6
$ws.Range("K22").Value = 4737.6
$ws.Range("M22").Value = -4438.6

$ws.Range("H32").Value = 1930.3833
$ws.Range("I32").Value = 1967.7413
$ws.Range("K32").Value = 1967.7413
$ws.Range("M32").Value = -1680.7413

$ws.Range("H41").Value = 22876.084
$ws.Range("I41").Value = 1566.3334
$ws.Range("K41").Value = 1566.3334
$ws.Range("M41").Value = -1152.3334

$ws.Range("H102").Value = 5645.1763
$ws.Range("I102").Value = 4228.6924
$ws.Range("K102").Value = 4228.6924
$ws.Range("M102").Value = -2606.6924

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1750
$ws.Range("I8").Value = 1250
$ws.Range("J8").Value = 2250
$ws.Range("K8").Value = 1250
$ws.Range("L8").Value = 2250
$ws.Range("M8").Value = -1110
$ws.Range("N8").Value = -2530

$ws.Range("H33").Value = 8888
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H130").Value = 36333
$ws.Range("I130").Value = 30000
$ws.Range("J130").Value = 39499.5
$ws.Range("K130").Value = 30000
$ws.Range("L130").Value = 39499.5
$ws.Range("M130").Value = -24980
$ws.Range("N130").Value = -49539.5

$ws.Range("H131").Value = 70199.39999999999
$ws.Range("J131").Value = 70199.39999999999
$ws.Range("L131").Value = 70199.39999999999
$ws.Range("N131").Value = -80279.39999999999

$ws.Range("H134").Value = 4677.0835
$ws.Range("I134").Value = 4677.0835
$ws.Range("K134").Value = 14031.2505
$ws.Range("M134").Value = -11496.2505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 10000
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("N11").Value = -10280

$ws.Range("H22").Value = 546.55554
$ws.Range("J22").Value = 785
$ws.Range("L22").Value = 785
$ws.Range("N22").Value = -1485

$ws.Range("H31").Value = 5842
$ws.Range("I31").Value = 6571.2856
$ws.Range("J31").Value = 5416.5835
$ws.Range("K31").Value = 6571.2856
$ws.Range("L31").Value = 5416.5835
$ws.Range("M31").Value = -6276.2856
$ws.Range("N31").Value = -6006.5835

$ws.Range("H34").Value = 5842
$ws.Range("I34").Value = 6571.2856
$ws.Range("J34").Value = 5416.5835
$ws.Range("K34").Value = 6571.2856
$ws.Range("L34").Value = 5416.5835
$ws.Range("M34").Value = -6369.2856
$ws.Range("N34").Value = -5820.5835

$ws.Range("H94").Value = 1049.8334
$ws.Range("I94").Value = 825
$ws.Range("J94").Value = 1162.25
$ws.Range("K94").Value = 825
$ws.Range("L94").Value = 1162.25
$ws.Range("M94").Value = -374
$ws.Range("N94").Value = -2064.25

$ws.Range("H99").Value = 5048.727
$ws.Range("I99").Value = 4410.75
$ws.Range("J99").Value = 6750
$ws.Range("K99").Value = 4410.75
$ws.Range("L99").Value = 6750
$ws.Range("M99").Value = -2912.75
$ws.Range("N99").Value = -9746

$ws.Range("H122").Value = 3288.4285
$ws.Range("I122").Value = 3272.111
$ws.Range("J122").Value = 3317.8
$ws.Range("K122").Value = 9816.332999999999
$ws.Range("L122").Value = 9953.400000000001
$ws.Range("M122").Value = -7366.332999999999
$ws.Range("N122").Value = -14853.4

$ws.Range("H126").Value = 5048.727
$ws.Range("I126").Value = 4410.75
$ws.Range("J126").Value = 6750
$ws.Range("K126").Value = 13232.25
$ws.Range("L126").Value = 20250
$ws.Range("M126").Value = -10762.25
$ws.Range("N126").Value = -25190

$ws.Range("H133").Value = 52679.5
$ws.Range("J133").Value = 54215.4
$ws.Range("L133").Value = 54215.4
$ws.Range("N133").Value = -59275.4

$ws.Range("H134").Value = 7115.5
$ws.Range("I134").Value = 6319.8
$ws.Range("J134").Value = 9502.6
$ws.Range("K134").Value = 18959.4
$ws.Range("L134").Value = 28507.8
$ws.Range("M134").Value = -16424.4
$ws.Range("N134").Value = -33577.8

$ws.Range("H141").Value = 38749.25
$ws.Range("I141").Value = 32498.5
$ws.Range("K141").Value = 32498.5
$ws.Range("M141").Value = -27318.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 68666.664
$ws.Range("I56").Value = 68666.664
$ws.Range("K56").Value = 68666.664
$ws.Range("M56").Value = -68136.664

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H113").Value = 760.6667
$ws.Range("J113").Value = 788
$ws.Range("L113").Value = 2364
$ws.Range("N113").Value = -6704

$ws.Range("H126").Value = 11499.833
$ws.Range("I126").Value = 7999.6665
$ws.Range("K126").Value = 23998.9995
$ws.Range("M126").Value = -19058.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 12580.875
$ws.Range("I2").Value = 92.57143000000001
$ws.Range("K2").Value = 92.57143000000001
$ws.Range("M2").Value = 20.42856999999999

$ws.Range("H69").Value = 33990
$ws.Range("J69").Value = 33990
$ws.Range("L69").Value = 33990
$ws.Range("N69").Value = -35488

$ws.Range("H72").Value = 33990
$ws.Range("J72").Value = 33990
$ws.Range("L72").Value = 101970
$ws.Range("N72").Value = -109458

$ws.Range("H122").Value = 2840.1724
$ws.Range("J122").Value = 2766.5
$ws.Range("L122").Value = 8299.5
$ws.Range("N122").Value = -13199.5

$ws.Range("H126").Value = 4400.483
$ws.Range("J126").Value = 4509.8887
$ws.Range("L126").Value = 13529.6661
$ws.Range("N126").Value = -18469.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3333.3333
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888

$ws.Range("H33").Value = 4000
$ws.Range("I33").Value = 4000
$ws.Range("K33").Value = 4000
$ws.Range("M33").Value = -3710

$ws.Range("H42").Value = 12500
$ws.Range("I42").Value = 10000
$ws.Range("K42").Value = 10000
$ws.Range("M42").Value = -9437

$ws.Range("H46").Value = 9890.451999999999
$ws.Range("I46").Value = 3654.3684
$ws.Range("K46").Value = 3654.3684
$ws.Range("M46").Value = -3466.3684

$ws.Range("H49").Value = 12500
$ws.Range("I49").Value = 10000
$ws.Range("K49").Value = 10000
$ws.Range("M49").Value = -9853

$ws.Range("H61").Value = 77734.69500000001
$ws.Range("I61").Value = 111769.445
$ws.Range("K61").Value = 111769.445
$ws.Range("M61").Value = -111567.445

$ws.Range("H100").Value = 6925.5713
$ws.Range("I100").Value = 6499.6665
$ws.Range("J100").Value = 7245
$ws.Range("K100").Value = 6499.6665
$ws.Range("L100").Value = 7245
$ws.Range("M100").Value = -5958.6665
$ws.Range("N100").Value = -8327

$ws.Range("H113").Value = 77734.69500000001
$ws.Range("I113").Value = 111769.445
$ws.Range("K113").Value = 111769.445
$ws.Range("M113").Value = -109599.445

$ws.Range("H122").Value = 4041.5557
$ws.Range("I122").Value = 3843.75
$ws.Range("K122").Value = 11531.25
$ws.Range("M122").Value = -9081.25

$ws.Range("H126").Value = 3333.3333
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9500
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H122").Value = 5250.0557
$ws.Range("I122").Value = 5201.0835
$ws.Range("J122").Value = 5348
$ws.Range("K122").Value = 15603.2505
$ws.Range("L122").Value = 16044
$ws.Range("M122").Value = -13153.2505
$ws.Range("N122").Value = -20944

$ws.Range("H132").Value = 3773.3438
$ws.Range("I132").Value = 3333.7856
$ws.Range("J132").Value = 6850.25
$ws.Range("K132").Value = 10001.3568
$ws.Range("L132").Value = 20550.75
$ws.Range("M132").Value = -7471.356800000001
$ws.Range("N132").Value = -25610.75
